$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New handoff entry: ee81f607-5d83-48c5-baf3-4679802b6fad.md
# Adds row 9 to the "Overview", "zh-cn" and "de-de" tables/sheets.
# ---------------------------------------------------------------------------

$fileId   = "ee81f607-5d83-48c5-baf3-4679802b6fad"
$mdName   = "$fileId.md"
$mdPath   = "e2e\$fileId.md"
$hoDate   = "2016-08-18 20:45:35"
$zhDate   = "2016-08-18 20:45:31"
$deDate   = "2016-08-18 20:45:35"
$zhXlf    = "$fileId.0a3cdf3b000333dba5f738625f463e86bbd6dbcc.zh-cn.xlf"
$deXlf    = "$fileId.0a3cdf3b000333dba5f738625f463e86bbd6dbcc.de-de.xlf"
$mdUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a1b9c3d5e7f14826394a5b6c7d8e9f0a1b2c3d4/e2e/$mdName"

# ---------------------------------------------------------------------------
# Sheet "Overview" (table3) -- columns:
#  A File Name | B Path And Name | C Extension | D Publish URL
#  E zh-cn | F de-de | G Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------

$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()

$wsOverview.Range("A9").Value = $mdName
$wsOverview.Range("B9").Value = $mdPath
$wsOverview.Range("C9").Value = ".md"
$wsOverview.Range("D9").Value = "'"
$wsOverview.Range("E9").Value = "Ready for handoff"
$wsOverview.Range("F9").Value = "Ready for handoff"
$wsOverview.Range("G9").Value = $hoDate
$wsOverview.Range("G9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B9"), $mdUrl, "", "", $mdPath) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (table1) -- columns:
#  A Source File Name | B File Extension | C Status | D Source Path
#  E Priority | F Content Duplicate | G Latest Handoff File | H Latest Handoff Datetime
#  I Latest Target File | J Latest Handback File | K Latest Handback DateTime
#  L Reference Tokens | M To be localized | N Dependency From | O Has metadata | P Error Detail
# ---------------------------------------------------------------------------

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$rowZhCn = $loZhCn.ListRows.Add()

$wsZhCn.Range("A9").Value = $mdName
$wsZhCn.Range("B9").Value = ".md"
$wsZhCn.Range("C9").Value = "Ready for handoff"
$wsZhCn.Range("D9").Value = "e2e"
$wsZhCn.Range("E9").Value = "ht"
$wsZhCn.Range("F9").Value = "'False"
$wsZhCn.Range("G9").Value = $zhXlf
$wsZhCn.Range("H9").Value = $zhDate
$wsZhCn.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I9").Value = "'"
$wsZhCn.Range("J9").Value = "'"
$wsZhCn.Range("K9").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L9").Value = "'"
$wsZhCn.Range("M9").Value = "'True"
$wsZhCn.Range("N9").Value = "'"
$wsZhCn.Range("O9").Value = "'False"
$wsZhCn.Range("P9").Value = "'"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A9"), $mdUrl, "", "", $mdName) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de" (table2) -- same column layout as zh-cn
# ---------------------------------------------------------------------------

$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$rowDeDe = $loDeDe.ListRows.Add()

$wsDeDe.Range("A9").Value = $mdName
$wsDeDe.Range("B9").Value = ".md"
$wsDeDe.Range("C9").Value = "Ready for handoff"
$wsDeDe.Range("D9").Value = "e2e"
$wsDeDe.Range("E9").Value = "ht"
$wsDeDe.Range("F9").Value = "'False"
$wsDeDe.Range("G9").Value = $deXlf
$wsDeDe.Range("H9").Value = $deDate
$wsDeDe.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I9").Value = "'"
$wsDeDe.Range("J9").Value = "'"
$wsDeDe.Range("K9").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L9").Value = "'"
$wsDeDe.Range("M9").Value = "'True"
$wsDeDe.Range("N9").Value = "'"
$wsDeDe.Range("O9").Value = "'False"
$wsDeDe.Range("P9").Value = "'"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A9"), $mdUrl, "", "", $mdName) | Out-Null
